$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rename the two "UserName..." pings to "Username..." (casing fix),
# keeping every other cell/value in the table untouched.
$ws.Range("A6").Value = "UsernamePing"
$ws.Range("A7").Value = "UsernameValidationPing"

# Move the active selection from A11 to A8.
$ws.Range("A8").Select() | Out-Null
